# Refresh market-board snapshot values (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the Leve profit-tracking sheets, per the scheduled data-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1229.3
$ws.Range("I4").Value = 1254.7778
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 1254.7778
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -1140.7778
$ws.Range("N4").Value = -1228
$ws.Range("H38").Value = 29
$ws.Range("I38").Value = 29
$ws.Range("K38").Value = 87
$ws.Range("M38").Value = 285
$ws.Range("H70").Value = 1507.1428
$ws.Range("I70").Value = 1162.625
$ws.Range("K70").Value = 3487.875
$ws.Range("M70").Value = -3217.875
$ws.Range("H73").Value = 1507.1428
$ws.Range("I73").Value = 1162.625
$ws.Range("K73").Value = 3487.875
$ws.Range("M73").Value = -2551.875
$ws.Range("H76").Value = 4743.207
$ws.Range("I76").Value = 3650.1
$ws.Range("J76").Value = 5318.5264
$ws.Range("K76").Value = 3650.1
$ws.Range("L76").Value = 5318.5264
$ws.Range("M76").Value = -3335.1
$ws.Range("N76").Value = -5948.5264
$ws.Range("H79").Value = 4743.207
$ws.Range("I79").Value = 3650.1
$ws.Range("J79").Value = 5318.5264
$ws.Range("K79").Value = 3650.1
$ws.Range("L79").Value = 5318.5264
$ws.Range("M79").Value = -2558.1
$ws.Range("N79").Value = -7502.5264
$ws.Range("H116").Value = 8151427.5
$ws.Range("I116").Value = 13583057
$ws.Range("K116").Value = 13583057
$ws.Range("M116").Value = -13579615
$ws.Range("H129").Value = 76924490
$ws.Range("I129").Value = 1279.2
$ws.Range("J129").Value = 333335200
$ws.Range("K129").Value = 3837.6
$ws.Range("L129").Value = 1000005600
$ws.Range("M129").Value = 1162.4
$ws.Range("N129").Value = -1000015600
$ws.Range("H132").Value = 3338011
$ws.Range("I132").Value = 4654.643
$ws.Range("K132").Value = 13963.929
$ws.Range("M132").Value = -11433.929
$ws.Range("H137").Value = 8303.5625
$ws.Range("J137").Value = 2402.6155
$ws.Range("L137").Value = 7207.8465
$ws.Range("N137").Value = -12307.8465
$ws.Range("H138").Value = 323160.38
$ws.Range("J138").Value = 4013.6155
$ws.Range("L138").Value = 12040.8465
$ws.Range("N138").Value = -22320.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 2990
$ws.Range("J30").Value = 2990
$ws.Range("L30").Value = 2990
$ws.Range("N30").Value = -3290
$ws.Range("H61").Value = 13144.9
$ws.Range("I61").Value = 22988.875
$ws.Range("K61").Value = 22988.875
$ws.Range("M61").Value = -22776.875
$ws.Range("H74").Value = 7714.7
$ws.Range("I74").Value = 9235.429
$ws.Range("K74").Value = 9235.429
$ws.Range("M74").Value = -8361.429
$ws.Range("H77").Value = 7714.7
$ws.Range("I77").Value = 9235.429
$ws.Range("K77").Value = 46177.145
$ws.Range("M77").Value = -41809.145
$ws.Range("H82").Value = 65000
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65722
$ws.Range("H85").Value = 65000
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67496
$ws.Range("H88").Value = 1253.0476
$ws.Range("I88").Value = 782.3333
$ws.Range("J88").Value = 1606.0834
$ws.Range("K88").Value = 782.3333
$ws.Range("L88").Value = 1606.0834
$ws.Range("M88").Value = -376.3333
$ws.Range("N88").Value = -2418.0834
$ws.Range("H91").Value = 1253.0476
$ws.Range("I91").Value = 782.3333
$ws.Range("J91").Value = 1606.0834
$ws.Range("K91").Value = 782.3333
$ws.Range("L91").Value = 1606.0834
$ws.Range("M91").Value = 621.6667
$ws.Range("N91").Value = -4414.0834
$ws.Range("H102").Value = 7422.8335
$ws.Range("I102").Value = 8870.5
$ws.Range("J102").Value = 5396.1
$ws.Range("K102").Value = 8870.5
$ws.Range("L102").Value = 5396.1
$ws.Range("M102").Value = -7248.5
$ws.Range("N102").Value = -8640.1
$ws.Range("H122").Value = 974362.9399999999
$ws.Range("I122").Value = 6133.0586
$ws.Range("K122").Value = 18399.1758
$ws.Range("M122").Value = -15949.1758
$ws.Range("H132").Value = 2916.8125
$ws.Range("I132").Value = 1458.3914
$ws.Range("J132").Value = 6643.8887
$ws.Range("K132").Value = 4375.174199999999
$ws.Range("L132").Value = 19931.6661
$ws.Range("M132").Value = -1845.174199999999
$ws.Range("N132").Value = -24991.6661
$ws.Range("H136").Value = 13144.9
$ws.Range("I136").Value = 22988.875
$ws.Range("K136").Value = 68966.625
$ws.Range("M136").Value = -66416.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7335.8
$ws.Range("I105").Value = 7357.6313
$ws.Range("K105").Value = 7357.6313
$ws.Range("M105").Value = -5610.6313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 22709.889
$ws.Range("I7").Value = 33614.832
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 33614.832
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -33501.832
$ws.Range("N7").Value = -1126
$ws.Range("H31").Value = 10200.125
$ws.Range("I31").Value = 12609.272
$ws.Range("J31").Value = 4900
$ws.Range("K31").Value = 12609.272
$ws.Range("L31").Value = 4900
$ws.Range("M31").Value = -12314.272
$ws.Range("N31").Value = -5490
$ws.Range("H34").Value = 10200.125
$ws.Range("I34").Value = 12609.272
$ws.Range("J34").Value = 4900
$ws.Range("K34").Value = 12609.272
$ws.Range("L34").Value = 4900
$ws.Range("M34").Value = -12407.272
$ws.Range("N34").Value = -5304
$ws.Range("H93").Value = 20782.834
$ws.Range("I93").Value = 9939.4
$ws.Range("K93").Value = 9939.4
$ws.Range("M93").Value = -8067.4
$ws.Range("H132").Value = 1830.6
$ws.Range("I132").Value = 1834.9231
$ws.Range("J132").Value = 1802.5
$ws.Range("K132").Value = 5504.7693
$ws.Range("L132").Value = 5407.5
$ws.Range("M132").Value = -2974.7693
$ws.Range("N132").Value = -10467.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 38175224
$ws.Range("I4").Value = 36831132
$ws.Range("K4").Value = 110493396
$ws.Range("M4").Value = -110493284
$ws.Range("H11").Value = 84.90000000000001
$ws.Range("I11").Value = 49.833332
$ws.Range("J11").Value = 137.5
$ws.Range("K11").Value = 149.499996
$ws.Range("L11").Value = 412.5
$ws.Range("M11").Value = -9.49999600000001
$ws.Range("N11").Value = -692.5
$ws.Range("H48").Value = 7047
$ws.Range("J48").Value = 7047
$ws.Range("L48").Value = 21141
$ws.Range("N48").Value = -21641

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3628
$ws.Range("I132").Value = 3763.5715
$ws.Range("J132").Value = 2950.1428
$ws.Range("K132").Value = 11290.7145
$ws.Range("L132").Value = 8850.428400000001
$ws.Range("M132").Value = -8760.7145
$ws.Range("N132").Value = -13910.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 648.6667
$ws.Range("J9").Value = 675
$ws.Range("L9").Value = 675
$ws.Range("N9").Value = -1123
$ws.Range("H22").Value = 8623.076999999999
$ws.Range("I22").Value = 17166.666
$ws.Range("K22").Value = 17166.666
$ws.Range("M22").Value = -16871.666
$ws.Range("H27").Value = 8623.076999999999
$ws.Range("I27").Value = 17166.666
$ws.Range("K27").Value = 17166.666
$ws.Range("M27").Value = -17059.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 13742.5
$ws.Range("J58").Value = 13500
$ws.Range("L58").Value = 13500
$ws.Range("N58").Value = -14116
$ws.Range("H132").Value = 4897.3335
$ws.Range("I132").Value = 4628.375
$ws.Range("K132").Value = 13885.125
$ws.Range("M132").Value = -11355.125
$ws.Range("H136").Value = 475184.44
$ws.Range("I136").Value = 537102.8
$ws.Range("J136").Value = 26276
$ws.Range("K136").Value = 1611308.4
$ws.Range("L136").Value = 78828
$ws.Range("M136").Value = -1608758.4
$ws.Range("N136").Value = -83928
